# Insert a new data row at row 36 (pushing existing rows 36.. down by one,
# e.g. old row 36 becomes row 37, ..., old row 131 becomes row 132) and
# populate the newly inserted row with the reported market data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 36..end down by one row.
$ws.Rows.Item(36).Insert()

# Fill in the new row 36 with the new daily price record.
$ws.Range("A36").Value = 3
$ws.Range("B36").Value = "Femacal de La Calera"
$ws.Range("C36").Value = "Coquimbo"
$ws.Range("D36").Value2 = 44708
$ws.Range("D36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E36").Value = 5
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100107
$ws.Range("H36").Value = "Otros"
$ws.Range("I36").Value = 100107011
$ws.Range("J36").Value = "Tuna"
$ws.Range("K36").Value = "Sin especificar"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 56
$ws.Range("N36").Value = 17000
$ws.Range("O36").Value = 17000
$ws.Range("P36").Value = 17000
$ws.Range("Q36").Value = "$/caja 20 kilos"
$ws.Range("R36").Value = "Provincia de Limarí"
$ws.Range("S36").Value = 850
$ws.Range("T36").Value = 20
